# Apply the crypto price/volume update described by the commit diff.
# Values in column D that are plain decimal numbers are entered with a
# leading apostrophe so Excel stores them as text (matching the sheet's
# existing inlineStr convention) instead of auto-converting to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "66.930.09"
$ws.Range("E2").Value = "  -0.70%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "3.457.67"
$ws.Range("E3").Value = "  -1.53%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  -0.04%  "

# --- Row 5: BNB ---
$ws.Range("D5").Value = "'592.85"
$ws.Range("E5").Value = "  -1.12%  "

# --- Row 6: Solana ---
$ws.Range("D6").Value = "'179.16"
$ws.Range("E6").Value = "  +1.65%  "

# --- Row 7: XRP ---
$ws.Range("E7").Value = "  +3.55%  "

# --- Row 8: USDC ---
$ws.Range("E8").Value = "  +0.04%  "

# --- Row 9: LidoStakedEther ---
$ws.Range("D9").Value = "3.453.47"
$ws.Range("E9").Value = "  -1.69%  "

# --- Row 10: Dogecoin ---
$ws.Range("E10").Value = "  +5.17%  "

# --- Row 11: Toncoin ---
$ws.Range("D11").Value = "'6.95"
$ws.Range("E11").Value = "  -3.32%  "

# --- Row 12: Cardano ---
$ws.Range("E12").Value = "  -0.48%  "

# --- Row 13: WrappedliquidstakedEther2.0 ---
$ws.Range("D13").Value = "4.053.34"
$ws.Range("E13").Value = "  -1.55%  "

# --- Row 14: Avalanche ---
$ws.Range("D14").Value = "'31.69"
$ws.Range("E14").Value = "  +2.95%  "

# --- Row 15: TRON ---
$ws.Range("E15").Value = "  -0.47%  "

# --- Row 16: WrappedBTC ---
$ws.Range("D16").Value = "66.855.79"
$ws.Range("E16").Value = "  -0.77%  "

# --- Row 17: ShibaInu ---
$ws.Range("E17").Value = "  -1.75%  "

# --- Row 18: WrappedEther ---
$ws.Range("D18").Value = "3.456.63"
$ws.Range("E18").Value = "  -1.08%  "

# --- Row 19: Polkadot ---
$ws.Range("E19").Value = "  -1.59%  "

# --- Row 20: Chainlink ---
$ws.Range("D20").Value = "'14.15"
$ws.Range("E20").Value = "  -3.28%  "

# --- Row 21: BitcoinCash ---
$ws.Range("D21").Value = "'388.41"
$ws.Range("E21").Value = "  -1.74%  "

# --- Row 22: Uniswap ---
$ws.Range("D22").Value = "'7.91"
$ws.Range("E22").Value = "  -1.35%  "

# --- Row 23: Dai ---
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.04%  "

# --- Row 24: LEO ---
$ws.Range("D24").Value = "'5.75"
$ws.Range("E24").Value = "  +1.22%  "

# --- Row 25: was Litecoin, now Polygon (rows 25/26 swapped) ---
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "'0.536"
$ws.Range("E25").Value = "  -0.64%  "

# --- Row 26: was Polygon, now Litecoin ---
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "'71.90"
$ws.Range("E26").Value = "  -2.29%  "

# --- Row 27: PEPE ---
$ws.Range("E27").Value = "  -1.07%  "

# --- Row 28: InternetComputer(DFINITY) ---
$ws.Range("D28").Value = "'10.30"
$ws.Range("E28").Value = "  +0.83%  "

# --- Row 29: Kaspa ---
$ws.Range("E29").Value = "  -3.71%  "

# --- Row 30: Binance-PegBSC-USD ---
$ws.Range("E30").Value = "  +0.55%  "

# --- Row 31: NEARProtocol ---
$ws.Range("D31").Value = "'6.15"
$ws.Range("E31").Value = "  -0.58%  "

# --- Row 32: Fetch.AI ---
$ws.Range("E32").Value = "  -2.37%  "

# --- Row 33: PancakeSwap ---
$ws.Range("E33").Value = "  -0.76%  "

# --- Row 34: EthereumClassic ---
$ws.Range("D34").Value = "'23.39"
$ws.Range("E34").Value = "  -1.36%  "

# --- Row 35: Aptos ---
$ws.Range("E35").Value = "  -1.57%  "

# --- Row 37: ImmutableX ---
$ws.Range("E37").Value = "  -3.95%  "

# --- Row 38: Monero ---
$ws.Range("D38").Value = "'163.16"
$ws.Range("E38").Value = "  -0.64%  "

# --- Row 39: Mantle ---
$ws.Range("D39").Value = "'0.876"
$ws.Range("E39").Value = "  -0.64%  "

# --- Row 40: dogwifhat ---
$ws.Range("D40").Value = "'2.80"
$ws.Range("E40").Value = "  +9.81%  "

# --- Row 41: Stacks ---
$ws.Range("D41").Value = "'1.87"
$ws.Range("E41").Value = "  -3.12%  "

# --- Row 42: RenderToken ---
$ws.Range("E42").Value = "  -4.62%  "

# --- Row 43: Filecoin ---
$ws.Range("E43").Value = "  -1.02%  "

# --- Row 44: EnergySwap ---
$ws.Range("D44").Value = "'26.13"
$ws.Range("E44").Value = "  -0.46%  "

# --- Row 45: Hedera ---
$ws.Range("D45").Value = "'0.0719"
$ws.Range("E45").Value = "  -2.19%  "

# --- Row 46: Maker ---
$ws.Range("D46").Value = "2.737.61"
$ws.Range("E46").Value = "  -2.46%  "

# --- Row 47: InjectiveProtocol ---
$ws.Range("D47").Value = "'26.04"
$ws.Range("E47").Value = "  -6.05%  "

# --- Row 48: OKB ---
$ws.Range("D48").Value = "'41.16"
$ws.Range("E48").Value = "  -3.30%  "

# --- Row 49: VeChain ---
$ws.Range("D49").Value = "'0.0298"
$ws.Range("E49").Value = "  -2.15%  "

# --- Row 50: Bittensor ---
$ws.Range("D50").Value = "'326.80"
$ws.Range("E50").Value = "  -4.74%  "

# --- Row 51: ONDO ---
$ws.Range("E51").Value = "  -4.11%  "
